$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the timesheet entries for rows 17-27 (previously blank placeholder rows)
$ws.Range("B17").Value = [DateTime]"2023-09-12"
$ws.Range("C17").Value = "DS-160-01"
$ws.Range("D17").Value = "80 minutes"
$ws.Range("E17").Value = "Lecture"

$ws.Range("B18").Value = [DateTime]"2023-09-12"
$ws.Range("C18").Value = "FILM-270-01"
$ws.Range("D18").Value = "120 minutes"
$ws.Range("E18").Value = "Homework"

$ws.Range("B19").Value = [DateTime]"2023-09-12"
$ws.Range("C19").Value = "FILM-270-01"
$ws.Range("D19").Value = "70 minutes"
$ws.Range("E19").Value = "Lecture"

$ws.Range("B20").Value = [DateTime]"2023-09-12"
$ws.Range("C20").Value = "MATH-205-03"
$ws.Range("D20").Value = "120 minutes"
$ws.Range("E20").Value = "Lecture"

$ws.Range("B21").Value = [DateTime]"2023-09-13"
$ws.Range("C21").Value = "CS-215-01"
$ws.Range("D21").Value = "50 minutes"
$ws.Range("E21").Value = "Lecture"

$ws.Range("B22").Value = [DateTime]"2023-09-13"
$ws.Range("C22").Value = "MUSI-111"
$ws.Range("D22").Value = "60 minutes"
$ws.Range("E22").Value = "Piano Lesson"

$ws.Range("B23").Value = [DateTime]"2023-09-13"
$ws.Range("C23").Value = "MATH-205-03"
$ws.Range("D23").Value = "75 minutes"
$ws.Range("E23").Value = "Lecture"

$ws.Range("B24").Value = [DateTime]"2023-09-17"
$ws.Range("C24").Value = "MATH-430-01"
$ws.Range("D24").Value = "60 minutes"
$ws.Range("E24").Value = "Homework"

$ws.Range("B25").Value = [DateTime]"2023-09-18"
$ws.Range("C25").Value = "CS-215-01"
$ws.Range("D25").Value = "50 minutes"
$ws.Range("E25").Value = "Lecture"

$ws.Range("B26").Value = [DateTime]"2023-09-18"
$ws.Range("C26").Value = "MATH-430-01"
$ws.Range("D26").Value = "75 minutes"
$ws.Range("E26").Value = "Lecture"

$ws.Range("B27").Value = [DateTime]"2023-09-18"
$ws.Range("C27").Value = "MUSE-133-01"
$ws.Range("D27").Value = "150 minutes"
$ws.Range("E27").Value = "Choir practice"

# Update the active selection to reflect where the user ended up (E28)
$ws.Range("E28").Select()
